$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("USB Oil 5 Early")

# Row 2
$ws.Range("I2").Value = 33
$ws.Range("J2").Value = "g"
$ws.Range("K2").Value = 1.5
$ws.Range("L2").Value = 37

# Row 3
$ws.Range("I3").Value = 34
$ws.Range("J3").Value = "t"
$ws.Range("K3").Value = 1.5
$ws.Range("L3").Value = 36

# Row 4
$ws.Range("I4").Value = 36
$ws.Range("J4").Value = "t"
$ws.Range("K4").Value = 1.5
$ws.Range("L4").Value = 33

# Row 5
$ws.Range("I5").Value = 33
$ws.Range("J5").Value = "t"
$ws.Range("K5").Value = 1.5
$ws.Range("L5").Value = 24

# Row 6
$ws.Range("I6").Value = 31
$ws.Range("J6").Value = "t"
$ws.Range("K6").Value = 1.5
$ws.Range("L6").Value = 28

# Row 7
$ws.Range("I7").Value = 31
$ws.Range("J7").Value = "g"
$ws.Range("K7").Value = 1.5
$ws.Range("L7").Value = 32

# Row 8
$ws.Range("I8").Value = 34
$ws.Range("J8").Value = "lt"
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 48
$ws.Range("P8").Value = "indet"

# Row 9
$ws.Range("I9").Value = 33
$ws.Range("J9").Value = "g"
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 36

# Row 10
$ws.Range("I10").Value = 32
$ws.Range("J10").Value = "g"
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 33

# Row 11
$ws.Range("I11").Value = 33
$ws.Range("J11").Value = "g"
$ws.Range("K11").Value = 1.5
$ws.Range("L11").Value = 33

# Row 12
$ws.Range("I12").Value = 31
$ws.Range("J12").Value = "g"
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 24

# Row 13
$ws.Range("I13").Value = 34
$ws.Range("J13").Value = "g"
$ws.Range("K13").Value = 1.5
$ws.Range("L13").Value = 35

# Row 14
$ws.Range("I14").Value = 35
$ws.Range("J14").Value = "t"
$ws.Range("K14").Value = 1.5
$ws.Range("L14").Value = 33

# Row 15
$ws.Range("I15").Value = 26
$ws.Range("J15").Value = "t"
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 31

# Row 16
$ws.Range("I16").Value = 34
$ws.Range("J16").Value = "g"
$ws.Range("K16").Value = 1.5
$ws.Range("L16").Value = 44

# Row 17
$ws.Range("I17").Value = 36
$ws.Range("J17").Value = "g"
$ws.Range("K17").Value = 1.5
$ws.Range("L17").Value = 29

# Row 18
$ws.Range("I18").Value = 33
$ws.Range("J18").Value = "t"
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 25

# Row 19
$ws.Range("I19").Value = 30
$ws.Range("J19").Value = "t"
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 21
$ws.Range("P19").Value = "thin"

# Row 20
$ws.Range("I20").Value = 33
$ws.Range("J20").Value = "g"
$ws.Range("K20").Value = 1.5
$ws.Range("L20").Value = 29

# Row 21
$ws.Range("I21").Value = 31
$ws.Range("J21").Value = "lt"
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 46
$ws.Range("P21").Value = "indet"

# Row 22
$ws.Range("I22").Value = "."
$ws.Range("J22").Value = "t"
$ws.Range("K22").Value = "."
$ws.Range("L22").Value = "."

# Row 23
$ws.Range("I23").Value = "."
$ws.Range("J23").Value = "g"
$ws.Range("K23").Value = "."
$ws.Range("L23").Value = "."

# Row 24
$ws.Range("I24").Value = "."
$ws.Range("J24").Value = "g"
$ws.Range("K24").Value = "."
$ws.Range("L24").Value = "."

# Row 25
$ws.Range("I25").Value = "."
$ws.Range("J25").Value = "t"
$ws.Range("K25").Value = "."
$ws.Range("L25").Value = "."

# Row 26
$ws.Range("I26").Value = "."
$ws.Range("J26").Value = "t"
$ws.Range("K26").Value = "."
$ws.Range("L26").Value = "."

# Row 27
$ws.Range("I27").Value = "."
$ws.Range("J27").Value = "g"
$ws.Range("K27").Value = "."
$ws.Range("L27").Value = "."

# Row 28
$ws.Range("I28").Value = "."
$ws.Range("J28").Value = "g"
$ws.Range("K28").Value = "."
$ws.Range("L28").Value = "."

# Row 29
$ws.Range("I29").Value = "."
$ws.Range("J29").Value = "g"
$ws.Range("K29").Value = "."
$ws.Range("L29").Value = "."

# Row 30
$ws.Range("I30").Value = "."
$ws.Range("J30").Value = "t"
$ws.Range("K30").Value = "."
$ws.Range("L30").Value = "."

# Row 31
$ws.Range("I31").Value = "."
$ws.Range("J31").Value = "lt"
$ws.Range("K31").Value = "."
$ws.Range("L31").Value = "."
$ws.Range("P31").Value = "indet"

# Refresh the view: scroll back to top, select the next empty column, and zoom in
$ws.Activate() | Out-Null
$ws.Range("O2").Select() | Out-Null
$excel.ActiveWindow.Zoom = 103
